$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new package entry in row 7, matching the pattern of existing rows
$ws.Range("A7").Value = "tespack"
$ws.Range("B7").Value = "tespack"
$ws.Range("C7").Value = "tespack"
$ws.Range("D7").Value = "1.0.0"
$ws.Range("E7").Value = "ContentPackage"

# F7 holds a date-like string ("2026-02-04") that must be stored as text,
# not auto-converted into a date serial number. Temporarily force text
# formatting while assigning the value, then restore the default style so
# the cell keeps using the workbook's normal/general formatting.
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2026-02-04"
$ws.Range("F7").Style = "Normal"
